$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the taxon/location data between rows 13 and 14 ---
$cols = @("A","B","E","F","G","H","P","Q","R","S")

$row13vals = @{}
$row14vals = @{}
foreach ($col in $cols) {
    $row13vals[$col] = $ws.Range("$col`13").Value2
    $row14vals[$col] = $ws.Range("$col`14").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col`13").Value = $row14vals[$col]
    $ws.Range("$col`14").Value = $row13vals[$col]
}

# --- 2. Round the Ost (Q) / Nord (R) coordinates to whole numbers for rows 12-14 ---
foreach ($row in 12..14) {
    $q = $ws.Range("Q$row").Value2
    $r = $ws.Range("R$row").Value2
    $ws.Range("Q$row").Value = [Math]::Round([double]$q, 0)
    $ws.Range("R$row").Value = [Math]::Round([double]$r, 0)
}

# --- 3. Clear the Starttid (Z) and Sluttid (AB) cells for rows 12-14 ---
foreach ($row in 12..14) {
    $ws.Range("Z$row").ClearContents()
    $ws.Range("AB$row").ClearContents()
}
